# Shubham Jain - Sprint-2 task tracking update
# Updates "Hours Burnt" (column G) for a handful of tasks and moves the
# sheet's active selection, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Hours Burnt (column G) for the affected tasks ---------------
# Column H ("Remaining Hours") holds =F-G formulas and recalculates
# automatically once G changes.
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 3
$ws.Range("G9").Value = 2
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 1

# --- Move the view/selection to match the saved workbook state ----------
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("G20").Select()
